$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text format first, to prevent Excel
# from silently converting numeric-looking strings (e.g. "0.3330")
# into real numbers and losing trailing zeros / exact formatting.
# NOTE: setting .NumberFormat on a multi-area (comma) Range only
# reliably affects the first area in this runtime, so loop instead.
$textCells = @("D2", "D3", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '22.467.05'
$ws.Range("D3").Value = '1.570.42'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").Value = '286.08'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").Value = '0.3649'
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").Value = '48.14'
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("D9").Value = '0.3330'
$ws.Range("E9").Value = '  -2.36%  '
$ws.Range("D10").Value = '1.128'
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("D11").Value = '0.07426'
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '20.86'
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("D14").Value = '5.969'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").Value = '6.912'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '1.569.15'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").Value = '0.00001106'
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("D18").Value = '87.98'
$ws.Range("E18").Value = '  -3.44%  '
$ws.Range("D19").Value = '0.06732'
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '6.378'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").Value = '16.42'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = '12.03'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").Value = '22.456.63'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = '2.633'
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").Value = '151.37'
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("D29").Value = '5.014'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '124.75'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").Value = '1.747.36'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '1.033'
$ws.Range("E32").Value = '  -4.48%  '
$ws.Range("D33").Value = '6.134'
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").Value = '1.993'
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").Value = '9.763'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").Value = '0.02424'
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("D38").Value = '0.2246'
$ws.Range("E38").Value = '  -2.67%  '
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("D40").Value = '5.414'
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").Value = '1.296'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").Value = '11.30'
$ws.Range("E42").Value = '  -0.76%  '
$ws.Range("D43").Value = '0.6285'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.90'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6096'
$ws.Range("E46").Value = '  +4.23%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.745'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.045'
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '124.06'
$ws.Range("E49").Value = '  -4.82%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.215'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.07226'
$ws.Range("E51").Value = '  -1.52%  '
